$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 90862
$ws.Cells.Item(2, 3).Value = 5539.9085
$ws.Cells.Item(2, 4).Value = 85322.09149999999

$ws.Cells.Item(3, 2).Value = 71987
$ws.Cells.Item(3, 3).Value = 5343.505999999999
$ws.Cells.Item(3, 4).Value = 66643.49400000001

$ws.Cells.Item(4, 2).Value = 69101
$ws.Cells.Item(4, 3).Value = 5250.640500000001
$ws.Cells.Item(4, 4).Value = 63850.3595

$ws.Cells.Item(5, 2).Value = 67174
$ws.Cells.Item(5, 3).Value = 5194.129499999999
$ws.Cells.Item(5, 4).Value = 61979.8705

$ws.Cells.Item(6, 2).Value = 67709
$ws.Cells.Item(6, 3).Value = 5143.397
$ws.Cells.Item(6, 4).Value = 62565.603

$ws.Cells.Item(7, 2).Value = 79680
$ws.Cells.Item(7, 3).Value = 5188.8305
$ws.Cells.Item(7, 4).Value = 74491.1695

$ws.Cells.Item(8, 2).Value = 66217
$ws.Cells.Item(8, 3).Value = 5628.136500000001
$ws.Cells.Item(8, 4).Value = 60588.8635

$ws.Cells.Item(9, 2).Value = 73755
$ws.Cells.Item(9, 3).Value = 6549.074000000001
$ws.Cells.Item(9, 4).Value = 67205.92599999999

$ws.Cells.Item(10, 2).Value = 82980
$ws.Cells.Item(10, 3).Value = 7809.641
$ws.Cells.Item(10, 4).Value = 75170.359

$ws.Cells.Item(11, 2).Value = 95273
$ws.Cells.Item(11, 3).Value = 12977.8705
$ws.Cells.Item(11, 4).Value = 82295.1295

$ws.Cells.Item(12, 2).Value = 97940
$ws.Cells.Item(12, 3).Value = 14879.494
$ws.Cells.Item(12, 4).Value = 83060.50599999999

$ws.Cells.Item(13, 2).Value = 95015
$ws.Cells.Item(13, 3).Value = 15334.284
$ws.Cells.Item(13, 4).Value = 79680.716

$ws.Cells.Item(14, 2).Value = 103218
$ws.Cells.Item(14, 3).Value = 15538.4565
$ws.Cells.Item(14, 4).Value = 87679.5435

$ws.Cells.Item(15, 2).Value = 100719
$ws.Cells.Item(15, 3).Value = 15473.9585
$ws.Cells.Item(15, 4).Value = 85245.04149999999

$ws.Cells.Item(16, 2).Value = 99562
$ws.Cells.Item(16, 3).Value = 15459.3285
$ws.Cells.Item(16, 4).Value = 84102.6715

$ws.Cells.Item(17, 2).Value = 98100
$ws.Cells.Item(17, 3).Value = 15733.5675
$ws.Cells.Item(17, 4).Value = 82366.4325

$ws.Cells.Item(18, 2).Value = 94912
$ws.Cells.Item(18, 3).Value = 15999.228
$ws.Cells.Item(18, 4).Value = 78912.772

$ws.Cells.Item(19, 2).Value = 93211
$ws.Cells.Item(19, 3).Value = 15687.07
$ws.Cells.Item(19, 4).Value = 77523.92999999999

$ws.Cells.Item(20, 2).Value = 95142
$ws.Cells.Item(20, 3).Value = 15318.933
$ws.Cells.Item(20, 4).Value = 79823.067

$ws.Cells.Item(21, 2).Value = 91234
$ws.Cells.Item(21, 3).Value = 13434.6415
$ws.Cells.Item(21, 4).Value = 77799.3585

$ws.Cells.Item(22, 2).Value = 88757
$ws.Cells.Item(22, 3).Value = 11609.815
$ws.Cells.Item(22, 4).Value = 77147.185

$ws.Cells.Item(23, 2).Value = 85054
$ws.Cells.Item(23, 3).Value = 9103.451000000001
$ws.Cells.Item(23, 4).Value = 75950.549

$ws.Cells.Item(24, 2).Value = 80176
$ws.Cells.Item(24, 3).Value = 6290.4415
$ws.Cells.Item(24, 4).Value = 73885.5585

$ws.Cells.Item(25, 2).Value = 76785
$ws.Cells.Item(25, 3).Value = 5217.8665
$ws.Cells.Item(25, 4).Value = 71567.1335

